$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 59408
$ws.Range("C10").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D10").Value = 388.17
$ws.Range("E10").Value = 463.78
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 2329.02
$ws.Range("B11").Value = 47438
$ws.Range("C11").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 401.81
$ws.Range("E11").Value = 480.05
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 803.62
$ws.Range("F48").Value = 23
$ws.Range("G48").Value = 4525.71
$ws.Range("F51").Value = 208
$ws.Range("G51").Value = 14568.32
$ws.Range("F58").Value = 36
$ws.Range("G58").Value = 1271.52
$ws.Range("F61").Value = 45
$ws.Range("G61").Value = 1138.05
$ws.Range("F63").Value = 25
$ws.Range("G63").Value = 938.5
$ws.Range("F78").Value = 67
$ws.Range("G78").Value = 1253.57
$ws.Range("B85").Value = 149841.23
$ws.Range("F136").Value = 51
$ws.Range("G136").Value = 1836
$ws.Range("F144").Value = 56
$ws.Range("G144").Value = 5112.24
$ws.Range("F146").Value = 29
$ws.Range("G146").Value = 6376.23
$ws.Range("F149").Value = 2
$ws.Range("G149").Value = 30.84
$ws.Range("F152").Value = 57
$ws.Range("G152").Value = 2418.51
$ws.Range("F153").Value = 35
$ws.Range("G153").Value = 855.05
$ws.Range("F154").Value = 43
$ws.Range("G154").Value = 4146.49
$ws.Range("B159").Value = 66077.56
$ws.Range("B193").Value = 53925
$ws.Range("B194").Value = 57756
$ws.Range("F196").Value = 65
$ws.Range("G196").Value = 2983.5
$ws.Range("B198").Value = 41722.54
$ws.Range("F226").Value = 48
$ws.Range("G226").Value = 2193.6
$ws.Range("B228").Value = 33938.28
$ws.Range("F293").Value = 24
$ws.Range("G293").Value = 16128.96
$ws.Range("B298").Value = 104306.49
$ws.Range("F425").Value = 9
$ws.Range("G425").Value = 1652.22
$ws.Range("B437").Value = 16664.75
$ws.Range("F457").Value = 11
$ws.Range("G457").Value = 233.75
$ws.Range("F472").Value = 19
$ws.Range("G472").Value = 817.1900000000001
$ws.Range("F480").Value = 53
$ws.Range("G480").Value = 6620.23
$ws.Range("B481").Value = 40024.11
$ws.Range("F495").Value = 12
$ws.Range("G495").Value = 1991.88
$ws.Range("B497").Value = 30764.77
$ws.Range("F561").Value = 860
$ws.Range("G561").Value = 11094
$ws.Range("B567").Value = 49116.29
$ws.Range("F591").Value = 304
$ws.Range("G591").Value = 11238.88
$ws.Range("F592").Value = 71
$ws.Range("G592").Value = 2222.3
$ws.Range("F596").Value = 171
$ws.Range("G596").Value = 4962.42
$ws.Range("B610").Value = 53209.13
$ws.Range("F618").Value = 599
$ws.Range("G618").Value = 6421.28
$ws.Range("B623").Value = 49151
$ws.Range("C623").Value = "NES-MAGGI Atta Noodles Masala 290g"
$ws.Range("D623").Value = 78.09999999999999
$ws.Range("E623").Value = 88.58
$ws.Range("F623").Value = 1
$ws.Range("G623").Value = 78.09999999999999
$ws.Range("B624").Value = 55667
$ws.Range("C624").Value = "NES-Maggi Atta Noodles Masala 290G"
$ws.Range("D624").Value = 85.76000000000001
$ws.Range("E624").Value = 97.25
$ws.Range("F624").Value = 55
$ws.Range("G624").Value = 4716.8
$ws.Range("B638").Value = 138778.34
$ws.Range("F660").Value = 17
$ws.Range("G660").Value = 2943.89
$ws.Range("F666").Value = 37
$ws.Range("G666").Value = 2562.99
$ws.Range("B667").Value = 24030
$ws.Range("F669").Value = 41
$ws.Range("G669").Value = 3248.43
$ws.Range("B688").Value = 82407.27
$ws.Range("F715").Value = 31
$ws.Range("G715").Value = 4047.05
$ws.Range("F719").Value = 79
$ws.Range("G719").Value = 2148.8
$ws.Range("B720").Value = 23024.09
$ws.Range("F818").Value = 189
$ws.Range("G818").Value = 9045.540000000001
$ws.Range("F820").Value = 49
$ws.Range("G820").Value = 3996.44
$ws.Range("F821").Value = 111
$ws.Range("G821").Value = 14774.1
$ws.Range("F825").Value = 18
$ws.Range("G825").Value = 671.04
$ws.Range("B837").Value = 187397.5
$ws.Range("F839").Value = 27
$ws.Range("G839").Value = 6735.96
$ws.Range("F843").Value = 58
$ws.Range("G843").Value = 6310.98
$ws.Range("F856").Value = 490
$ws.Range("G856").Value = 14758.8
$ws.Range("F858").Value = 35
$ws.Range("G858").Value = 13223.35
$ws.Range("F862").Value = 7
$ws.Range("G862").Value = 330.47
$ws.Range("F863").Value = 134
$ws.Range("G863").Value = 4262.54
$ws.Range("B867").Value = 191220.56
$ws.Range("B923").Value = 2427909.83
$ws.Range("B924").Value = 2427909.83
